# Apply "fix errors && add features && removed outdated" edits.
#
# 1. Sheet1: move the active selection from S20 to C15 (this also drops the
#    stale topLeftCell="M1" scroll position that was left over from a
#    previous view).
# 2. Sheet1 rows 3-10: wire up the second/third/fourth "start time" cells
#    (N, T, Z) so that they mirror the corresponding "end time" cell from the
#    previous slot (I, O, U respectively) via a formula instead of being
#    hard-typed / blank.
# 3. Sheet1 rows 3-10: the AG column (a straight "minutes worked" calc) now
#    also folds in the AF column (overtime/bonus minutes) so the daily total
#    lines up correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 3; $row -le 10; $row++) {
    $ws.Range("N$row").Formula  = "=I$row"
    $ws.Range("T$row").Formula  = "=O$row"
    $ws.Range("Z$row").Formula  = "=U$row"
    $ws.Range("AG$row").Formula = "=(MINUTE(D$row)*60+SECOND(D$row))*10/60+AF$row"
}

# Update the saved selection / scroll state on Sheet1.
[void]$ws.Activate()
[void]$ws.Range("C15").Select()
